{"js": "// 1. Merge the two runs \"{structure2\" + \"}\" (which are split apart by the\n//    \"_GoBack\" bookmark) back into a single run \"{structure2}\", removing\n//    the bookmark that sat between them in the process.\nconst mergeResults = context.document.body.search(\"{structure2}\", { matchCase: true, matchWildcards: false });\nmergeResults.load(\"items\");\nawait context.sync();\n\nif (mergeResults.items.length > 0) {\n  const mergeRange = mergeResults.items[0];\n  // Re-inserting the same visible text over the matched range collapses the\n  // two runs (and drops the now-empty bookmark that separated them) into a\n  // single run, while keeping the formatting of the first run.\n  mergeRange.insertText(\"{structure2}\", \"Replace\");\n  await context.sync();\n}\n\n// 2. Re-create the \"_GoBack\" bookmark at the end of the (empty) paragraph\n//    that carries the two-column section break right after the \"F\u0259hl\u0259\"\n//    paragraph.\nconst anchorResults = context.document.body.search(\"F\u0259hl\u0259\", { matchCase: true, matchWildcards: false });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length > 0) {\n  const anchorParagraph = anchorResults.items[0].paragraphs.getFirst();\n  const targetParagraph = anchorParagraph.getNext();\n  const endOfParagraph = targetParagraph.getRange(\"End\");\n  endOfParagraph.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Merge the two runs \"{structure2\" + \"}\" (they are split apart by the\n#    \"_GoBack\" bookmark) back into a single run \"{structure2}\", removing the\n#    bookmark that sat between them in the process.\n#\n#    Find/Execute locates the (bookmark-spanning) text as one contiguous\n#    range. Re-assigning Range.Text to the exact same string is a no-op in\n#    this engine, so we first set it to a distinct placeholder (forcing the\n#    underlying runs/bookmark to collapse) and then set it back to the\n#    original text.\n$mergeRange = $d.Content\n$mergeFind = $mergeRange.Find\n$mergeFind.ClearFormatting()\n$mergeFind.MatchWildcards = $false\n$mergeFind.Text = \"{structure2}\"\n$mergeFound = $mergeFind.Execute()\nif ($mergeFound) {\n    $mergeRange.Text = \"{structure2}__tmp__\"\n    $mergeRange.Text = \"{structure2}\"\n}\n\n# 2. Re-create the \"_GoBack\" bookmark at the end of the (empty) paragraph\n#    that carries the two-column section break, i.e. the paragraph right\n#    after the one containing \"F\u0259hl\u0259\".\n$anchorRange = $d.Content\n$anchorFind = $anchorRange.Find\n$anchorFind.ClearFormatting()\n$anchorFind.MatchWildcards = $false\n$anchorFind.Text = \"F\u0259hl\u0259\"\n$anchorFound = $anchorFind.Execute()\nif ($anchorFound) {\n    $anchorParagraph = $anchorRange.Paragraphs(1)\n    $targetParagraph = $anchorParagraph.Next()\n    $bookmarkRange = $targetParagraph.Range.Duplicate()\n    $bookmarkRange.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n}\n"}
